# "Generate Report for Handback"
#
# Updates the localization-status report to reflect a completed handback:
#  - Status text flips from "Ready for handoff" to "Handed back: in sync with en-US"
#  - The per-language tables (zh-cn, de-de) get their "Latest Target File" (hyperlinked
#    to the source .md on GitHub) and "Latest Handback File" (xlf name) columns filled in
#  - de-de additionally gets its "Latest Handback DateTime" stamped
#  - A few columns are widened so the new long file names are readable

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

$mdUrlA = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0234baa273f242caedbdbcba59a55d1591876c24/e2e/a6145d5c-2579-48fa-b66d-8afb04c48a5e.md"
$mdUrlC = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0234baa273f242caedbdbcba59a55d1591876c24/e2e/c5fcb705-c6e3-4f9d-a60f-23cf3b232c2e.md"
$mdNameA = "a6145d5c-2579-48fa-b66d-8afb04c48a5e.md"
$mdNameC = "c5fcb705-c6e3-4f9d-a60f-23cf3b232c2e.md"

function Set-HyperlinkLook($range) {
    # Match the look of the workbook's existing custom "HyperLink" cell style
    # (underline + cornflower-blue font) that rows A2/A3 already use.
    $range.Font.Underline = 2
    $range.Font.Color = 15570276
}

# ---- Overview sheet: just the status text (shared with the detail sheets) ----
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("E2").Value = $statusText
$ov.Range("F2").Value = $statusText
$ov.Range("E3").Value = $statusText
$ov.Range("F3").Value = $statusText
$ov.Columns.Item(5).ColumnWidth = 29.17
$ov.Columns.Item(6).ColumnWidth = 29.17

# ---- zh-cn detail sheet ----
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("C2").Value = $statusText
$zh.Range("C3").Value = $statusText

$zh.Range("I2").Value = $mdNameA
$zh.Hyperlinks.Add($zh.Range("I2"), $mdUrlA, [Type]::Missing, [Type]::Missing, $mdNameA)
Set-HyperlinkLook $zh.Range("I2")

$zh.Range("I3").Value = $mdNameC
$zh.Hyperlinks.Add($zh.Range("I3"), $mdUrlC, [Type]::Missing, [Type]::Missing, $mdNameC)
Set-HyperlinkLook $zh.Range("I3")

$zh.Range("J2").Value = "a6145d5c-2579-48fa-b66d-8afb04c48a5e.5a13e7ca6a7e333eacbc6bb7219076c1231dc501.zh-cn.xlf"
$zh.Range("J3").Value = "c5fcb705-c6e3-4f9d-a60f-23cf3b232c2e.30c397a142a363042ffa252376cba969e73ff706.zh-cn.xlf"

$zh.Range("K2").Value = "2016-09-06 06:34:00"
$zh.Range("K3").Value = "2016-09-06 06:34:00"

$zh.Columns.Item(3).ColumnWidth = 29.17
$zh.Columns.Item(9).ColumnWidth = 39.17
$zh.Columns.Item(10).ColumnWidth = 39.17

# ---- de-de detail sheet ----
$de = $wb.Worksheets.Item("de-de")

$de.Range("C2").Value = $statusText
$de.Range("C3").Value = $statusText

$de.Range("I2").Value = $mdNameA
$de.Hyperlinks.Add($de.Range("I2"), $mdUrlA, [Type]::Missing, [Type]::Missing, $mdNameA)
Set-HyperlinkLook $de.Range("I2")

$de.Range("I3").Value = $mdNameC
$de.Hyperlinks.Add($de.Range("I3"), $mdUrlC, [Type]::Missing, [Type]::Missing, $mdNameC)
Set-HyperlinkLook $de.Range("I3")

$de.Range("J2").Value = "a6145d5c-2579-48fa-b66d-8afb04c48a5e.5a13e7ca6a7e333eacbc6bb7219076c1231dc501.de-de.xlf"
$de.Range("J3").Value = "c5fcb705-c6e3-4f9d-a60f-23cf3b232c2e.30c397a142a363042ffa252376cba969e73ff706.de-de.xlf"

$de.Range("K2").Value = "2016-09-06 06:34:17"
$de.Range("K3").Value = "2016-09-06 06:34:17"

$de.Columns.Item(3).ColumnWidth = 29.17
$de.Columns.Item(9).ColumnWidth = 39.17
$de.Columns.Item(10).ColumnWidth = 39.17
